$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22 - this shifts the existing rows 22-37
# down to 23-38, preserving all of their data/formatting.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record.
$ws.Range("A22").Value = 8
$ws.Range("B22").Value = "Terminal La Palmera de La Serena"
$ws.Range("C22").Value = "Coquimbo"
$ws.Range("D22").Value = 44873
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100101
$ws.Range("H22").Value = "Berries"
$ws.Range("I22").Value = 100101001
$ws.Range("J22").Value = "Arándano (blue)"
$ws.Range("K22").Value = "Sin especificar"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 9000
$ws.Range("O22").Value = 10000
$ws.Range("P22").Value = 9500
$ws.Range("Q22").Value = "$/bandeja 2 kilos"
$ws.Range("R22").Value = "Provincia de Limarí"
$ws.Range("S22").Value = 4750
$ws.Range("T22").Value = 2
